# Updates cryptos list data (prices / 1h volume % change) on the active sheet.
# All Price (D) and Volume(1h) (E) values are stored as plain text in the sheet,
# so for cells whose new value would otherwise be auto-parsed by Excel as a
# number (single-decimal values like "214.27"), we prefix with a literal
# leading apostrophe ('' inside a single-quoted PowerShell string produces a
# single literal apostrophe) to force Excel to keep/store them as text, just
# like the original values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value  = '27.058.79'
$ws.Range('E2').Value  = '  -0.59%  '

$ws.Range('D3').Value  = '1.621.55'
$ws.Range('E3').Value  = '  -1.39%  '

$ws.Range('E4').Value  = '  -0.13%  '

$ws.Range('D5').Value  = '''214.27'
$ws.Range('E5').Value  = '  -1.36%  '

$ws.Range('D6').Value  = '''0.519'
$ws.Range('E6').Value  = '  +0.88%  '

$ws.Range('E7').Value  = '  -0.11%  '

$ws.Range('E8').Value  = '  -1.55%  '

$ws.Range('E9').Value  = '  -0.25%  '

$ws.Range('D10').Value = '''20.31'
$ws.Range('E10').Value = '  +1.78%  '

$ws.Range('E11').Value = '  -0.12%  '

$ws.Range('D12').Value = '1.624.11'
$ws.Range('E12').Value = '  -1.19%  '

$ws.Range('E13').Value = '  -0.51%  '

$ws.Range('D14').Value = '''0.541'
$ws.Range('E14').Value = '  -0.57%  '

$ws.Range('D15').Value = '27.029.40'
$ws.Range('E15').Value = '  -0.66%  '

$ws.Range('D16').Value = '''64.40'
$ws.Range('E16').Value = '  -4.46%  '

$ws.Range('D17').Value = '0.0₃0743'
$ws.Range('E17').Value = '  +0.64%  '

$ws.Range('D18').Value = '''215.89'
$ws.Range('E18').Value = '  -1.38%  '

$ws.Range('E19').Value = '  -0.05%  '

$ws.Range('E20').Value = '  +0.70%  '

$ws.Range('E21').Value = '  -0.71%  '

$ws.Range('D22').Value = '''2.41'
$ws.Range('E22').Value = '  -5.66%  '

$ws.Range('E23').Value = '  -1.68%  '

$ws.Range('D24').Value = '''147.11'
$ws.Range('E24').Value = '  -0.59%  '

$ws.Range('E25').Value = '  -0.20%  '

$ws.Range('E26').Value = '  -3.49%  '

$ws.Range('E27').Value = '  -0.47%  '

$ws.Range('E28').Value = '  -1.08%  '

$ws.Range('D29').Value = '''0.0504'
$ws.Range('E29').Value = '  -0.74%  '

$ws.Range('E30').Value = '  -1.03%  '

$ws.Range('E31').Value = '  -1.25%  '

$ws.Range('E32').Value = '  -1.52%  '

$ws.Range('D33').Value = '1.340.00'
$ws.Range('E33').Value = '  +6.36%  '

$ws.Range('E34').Value = '  -0.34%  '

$ws.Range('E35').Value = '  -0.44%  '

$ws.Range('E36').Value = '  -0.98%  '

$ws.Range('D37').Value = '''0.547'
$ws.Range('E37').Value = '  +0.41%  '

$ws.Range('D38').Value = '''0.852'
$ws.Range('E38').Value = '  +0.15%  '

$ws.Range('E39').Value = '  -0.14%  '

$ws.Range('D40').Value = '''0.801'
$ws.Range('E40').Value = '  -0.91%  '

$ws.Range('E41').Value = '  -0.12%  '

$ws.Range('D42').Value = '''65.30'
$ws.Range('E42').Value = '  +5.68%  '

$ws.Range('D43').Value = '1.758.23'
$ws.Range('E43').Value = '  -1.46%  '

$ws.Range('E44').Value = '  -1.89%  '

$ws.Range('D45').Value = '''90.37'
$ws.Range('E45').Value = '  -1.38%  '

# Rows 46 & 47 swap content (WEMIXToken / RenderToken trade ranking places)
$ws.Range('B46').Value = 'WEMIXToken'
$ws.Range('C46').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D46').Value = '''0.862'
$ws.Range('E46').Value = '  +29.87%  '

$ws.Range('B47').Value = 'RenderToken'
$ws.Range('C47').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D47').Value = '''1.61'
$ws.Range('E47').Value = '  +0.85%  '

$ws.Range('E48').Value = '  +3.63%  '

$ws.Range('E49').Value = '  -0.35%  '

$ws.Range('E50').Value = '  +2.10%  '

$ws.Range('D51').Value = '''7.56'
$ws.Range('E51').Value = '  -0.95%  '
